# Updates the crypto market snapshot on Sheet1 (columns B:E, rows 2-51)
# with freshly scraped values, matching the upstream "cryptos list" update.
#
# Columns:
#   B = Coin name, C = coinranking.com link, D = Price, E = Volume(1h) %
#
# Several "Price" values look like plain numbers (e.g. "45.56", "0.0000266")
# even though the source column stores them as text. Left alone, Excel's COM
# layer auto-converts such literals into real numeric cells (losing
# significant trailing zeros, or flipping tiny decimals into scientific
# notation). To avoid that, those specific cells are briefly switched to a
# Text ("@") number format before the value is assigned, then restored to the
# workbook's default "Normal" style so no stray formatting is left behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param([string]$CellRef, [string]$Text, [switch]$ForceText)

    $range = $ws.Range($CellRef)
    if ($ForceText) {
        $range.NumberFormat = "@"
        $range.Value = $Text
        $range.Style = "Normal"
    } else {
        $range.Value = $Text
    }
}


$updates = @(
    @{ Cell = "D2"; Value = '92.151.73'; ForceText = $false }
    @{ Cell = "E2"; Value = '  +5.02%  '; ForceText = $false }
    @{ Cell = "D3"; Value = '3.274.49'; ForceText = $false }
    @{ Cell = "E3"; Value = '  +0.28%  '; ForceText = $false }
    @{ Cell = "E4"; Value = '  +0.14%  '; ForceText = $false }
    @{ Cell = "D5"; Value = '216.50'; ForceText = $true }
    @{ Cell = "E5"; Value = '  +2.09%  '; ForceText = $false }
    @{ Cell = "D6"; Value = '628.30'; ForceText = $true }
    @{ Cell = "E6"; Value = '  +0.07%  '; ForceText = $false }
    @{ Cell = "D7"; Value = '0.413'; ForceText = $true }
    @{ Cell = "E7"; Value = '  +8.88%  '; ForceText = $false }
    @{ Cell = "D8"; Value = '0.721'; ForceText = $true }
    @{ Cell = "E8"; Value = '  +4.18%  '; ForceText = $false }
    @{ Cell = "E9"; Value = '  +0.02%  '; ForceText = $false }
    @{ Cell = "D10"; Value = '3.267.03'; ForceText = $false }
    @{ Cell = "E10"; Value = '  +0.26%  '; ForceText = $false }
    @{ Cell = "D11"; Value = '0.587'; ForceText = $true }
    @{ Cell = "E11"; Value = '  +1.46%  '; ForceText = $false }
    @{ Cell = "B12"; Value = 'ShibaInu'; ForceText = $false }
    @{ Cell = "C12"; Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'; ForceText = $false }
    @{ Cell = "D12"; Value = '0.0000266'; ForceText = $true }
    @{ Cell = "E12"; Value = '  +2.53%  '; ForceText = $false }
    @{ Cell = "B13"; Value = 'TRON'; ForceText = $false }
    @{ Cell = "C13"; Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'; ForceText = $false }
    @{ Cell = "D13"; Value = '0.180'; ForceText = $true }
    @{ Cell = "E13"; Value = '  -1.47%  '; ForceText = $false }
    @{ Cell = "D14"; Value = '34.24'; ForceText = $true }
    @{ Cell = "E14"; Value = '  +0.25%  '; ForceText = $false }
    @{ Cell = "D15"; Value = '3.878.20'; ForceText = $false }
    @{ Cell = "E15"; Value = '  +0.34%  '; ForceText = $false }
    @{ Cell = "D16"; Value = '92.030.61'; ForceText = $false }
    @{ Cell = "E16"; Value = '  +5.47%  '; ForceText = $false }
    @{ Cell = "D17"; Value = '5.34'; ForceText = $true }
    @{ Cell = "E17"; Value = '  +0.11%  '; ForceText = $false }
    @{ Cell = "D18"; Value = '3.256.71'; ForceText = $false }
    @{ Cell = "E18"; Value = '  -0.01%  '; ForceText = $false }
    @{ Cell = "D19"; Value = '3.33'; ForceText = $true }
    @{ Cell = "E19"; Value = '  +6.23%  '; ForceText = $false }
    @{ Cell = "D20"; Value = '14.08'; ForceText = $true }
    @{ Cell = "E20"; Value = '  +0.15%  '; ForceText = $false }
    @{ Cell = "D21"; Value = '439.34'; ForceText = $true }
    @{ Cell = "E21"; Value = '  +1.17%  '; ForceText = $false }
    @{ Cell = "B22"; Value = 'PEPE'; ForceText = $false }
    @{ Cell = "C22"; Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'; ForceText = $false }
    @{ Cell = "D22"; Value = '0.0000194'; ForceText = $true }
    @{ Cell = "E22"; Value = '  +49.03%  '; ForceText = $false }
    @{ Cell = "B23"; Value = 'Uniswap'; ForceText = $false }
    @{ Cell = "C23"; Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'; ForceText = $false }
    @{ Cell = "D23"; Value = '8.91'; ForceText = $true }
    @{ Cell = "E23"; Value = '  -0.01%  '; ForceText = $false }
    @{ Cell = "E24"; Value = '  -1.14%  '; ForceText = $false }
    @{ Cell = "D25"; Value = '5.37'; ForceText = $true }
    @{ Cell = "E25"; Value = '  +4.64%  '; ForceText = $false }
    @{ Cell = "D26"; Value = '12.35'; ForceText = $true }
    @{ Cell = "E26"; Value = '  -1.09%  '; ForceText = $false }
    @{ Cell = "D27"; Value = '3.460.81'; ForceText = $false }
    @{ Cell = "E27"; Value = '  +2.37%  '; ForceText = $false }
    @{ Cell = "D28"; Value = '77.12'; ForceText = $true }
    @{ Cell = "E28"; Value = '  +0.74%  '; ForceText = $false }
    @{ Cell = "E29"; Value = '  +0.02%  '; ForceText = $false }
    @{ Cell = "D30"; Value = '0.181'; ForceText = $true }
    @{ Cell = "E30"; Value = '  +0.37%  '; ForceText = $false }
    @{ Cell = "E31"; Value = '  +0.10%  '; ForceText = $false }
    @{ Cell = "D32"; Value = '8.77'; ForceText = $true }
    @{ Cell = "E32"; Value = '  -0.01%  '; ForceText = $false }
    @{ Cell = "D33"; Value = '552.67'; ForceText = $true }
    @{ Cell = "E33"; Value = '  +0.69%  '; ForceText = $false }
    @{ Cell = "D34"; Value = '7.12'; ForceText = $true }
    @{ Cell = "E34"; Value = '  +2.17%  '; ForceText = $false }
    @{ Cell = "E35"; Value = '  -1.53%  '; ForceText = $false }
    @{ Cell = "E36"; Value = '  +22.98%  '; ForceText = $false }
    @{ Cell = "E37"; Value = '  -8.43%  '; ForceText = $false }
    @{ Cell = "D38"; Value = '22.66'; ForceText = $true }
    @{ Cell = "E38"; Value = '  +0.75%  '; ForceText = $false }
    @{ Cell = "E39"; Value = '  +3.56%  '; ForceText = $false }
    @{ Cell = "E40"; Value = '  -4.65%  '; ForceText = $false }
    @{ Cell = "E42"; Value = '  +0.08%  '; ForceText = $false }
    @{ Cell = "E43"; Value = '  -0.31%  '; ForceText = $false }
    @{ Cell = "E44"; Value = '  -0.01%  '; ForceText = $false }
    @{ Cell = "D45"; Value = '150.66'; ForceText = $true }
    @{ Cell = "E45"; Value = '  -2.77%  '; ForceText = $false }
    @{ Cell = "B46"; Value = 'OKB'; ForceText = $false }
    @{ Cell = "C46"; Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'; ForceText = $false }
    @{ Cell = "D46"; Value = '45.56'; ForceText = $true }
    @{ Cell = "E46"; Value = '  +1.47%  '; ForceText = $false }
    @{ Cell = "B47"; Value = 'Aave'; ForceText = $false }
    @{ Cell = "C47"; Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'; ForceText = $false }
    @{ Cell = "D47"; Value = '180.09'; ForceText = $true }
    @{ Cell = "E47"; Value = '  +0.19%  '; ForceText = $false }
    @{ Cell = "D48"; Value = '0.130'; ForceText = $true }
    @{ Cell = "E48"; Value = '  +5.32%  '; ForceText = $false }
    @{ Cell = "D49"; Value = '1.27'; ForceText = $true }
    @{ Cell = "E49"; Value = '  -1.47%  '; ForceText = $false }
    @{ Cell = "B50"; Value = 'Filecoin'; ForceText = $false }
    @{ Cell = "C50"; Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'; ForceText = $false }
    @{ Cell = "D50"; Value = '4.23'; ForceText = $true }
    @{ Cell = "E50"; Value = '  +0.03%  '; ForceText = $false }
    @{ Cell = "B51"; Value = 'ARBITRUM'; ForceText = $false }
    @{ Cell = "C51"; Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'; ForceText = $false }
    @{ Cell = "D51"; Value = '0.635'; ForceText = $true }
    @{ Cell = "E51"; Value = '  +1.85%  '; ForceText = $false }
)

foreach ($update in $updates) {
    Set-TextValue $update.Cell $update.Value $update.ForceText
}

